# Update as per reviews:
# Add two new review comment rows to the "Code" worksheet, and make the
# "Code" sheet the active tab (it was previously "Design Document").

$wb = $excel.ActiveWorkbook

$codeSheet = $wb.Worksheets.Item("Code")

# New review comment row: negative-number transfer message
$codeSheet.Range("A5").Value = "Version 3"
$codeSheet.Range("B5").Value = "add specific messages to transfer`nmoney, if user entered negative `nnumber"
$codeSheet.Range("B5").WrapText = $true
$codeSheet.Range("C5").Value = "Mohamed Hassan"
$codeSheet.Range("D5").Value = "Closed"
$codeSheet.Rows.Item(5).RowHeight = 75

# New review comment row: zero-amount transfer message
$codeSheet.Range("A6").Value = "Version 3"
$codeSheet.Range("B6").Value = "add specific messages to transfer`nmoney, if user entered zeros"
$codeSheet.Range("B6").WrapText = $true
$codeSheet.Range("C6").Value = "Mohamed Hassan"
$codeSheet.Range("D6").Value = "Closed"
$codeSheet.Rows.Item(6).RowHeight = 45

# Move the selected/active cell + active tab from "Design Document" to "Code"
$codeSheet.Range("F5").Select()
$codeSheet.Activate()
